$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected ICDC Breed 1-14 FilesTab Neo4j script:
# remove the `File Type` and `Breed` columns from the FilesTab query (row 4, column B).
$newFilesTabQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Boston Terrier']   
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFilesTabQuery

# The query text got two lines shorter, so the row shrinks accordingly.
$ws.Rows(4).RowHeight = 217.5

# Move the selection/active cell to the cell that was edited.
[void]$ws.Range("B4").Select()
